$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped from the coinranking.com crypto-price refresh
# (Coin/Link/Price/Volume(1h) table on Sheet1). Each entry is the new
# cell text as it appears in the source spreadsheet.
$updates = @(
    @{ Cell = "D2"; Value = "29.716.53" },
    @{ Cell = "E2"; Value = "  +8.54%  " },
    @{ Cell = "D3"; Value = "1.946.51" },
    @{ Cell = "E3"; Value = "  +7.09%  " },
    @{ Cell = "E4"; Value = "  -0.40%  " },
    @{ Cell = "D5"; Value = "341.74" },
    @{ Cell = "E5"; Value = "  +3.01%  " },
    @{ Cell = "E6"; Value = "  -0.29%  " },
    @{ Cell = "D7"; Value = "0.4781" },
    @{ Cell = "E7"; Value = "  +4.70%  " },
    @{ Cell = "D8"; Value = "0.4136" },
    @{ Cell = "E8"; Value = "  +8.66%  " },
    @{ Cell = "E9"; Value = "  +5.34%  " },
    @{ Cell = "D10"; Value = "0.08254" },
    @{ Cell = "E10"; Value = "  +5.16%  " },
    @{ Cell = "D11"; Value = "1.041" },
    @{ Cell = "E11"; Value = "  +8.67%  " },
    @{ Cell = "D12"; Value = "22.67" },
    @{ Cell = "E12"; Value = "  +8.13%  " },
    @{ Cell = "D13"; Value = "1.930.61" },
    @{ Cell = "E13"; Value = "  +5.81%  " },
    @{ Cell = "D14"; Value = "6.183" },
    @{ Cell = "D15"; Value = "7.427" },
    @{ Cell = "E15"; Value = "  +5.26%  " },
    @{ Cell = "D16"; Value = "92.46" },
    @{ Cell = "E16"; Value = "  +3.45%  " },
    @{ Cell = "D17"; Value = "1.001" },
    @{ Cell = "E17"; Value = "  -0.26%  " },
    @{ Cell = "D18"; Value = "0.00001062" },
    @{ Cell = "E18"; Value = "  +4.23%  " },
    @{ Cell = "D19"; Value = "0.06682" },
    @{ Cell = "E19"; Value = "  +1.41%  " },
    @{ Cell = "D20"; Value = "18.09" },
    @{ Cell = "E20"; Value = "  +5.76%  " },
    @{ Cell = "E21"; Value = "  -0.24%  " },
    @{ Cell = "D22"; Value = "29.678.18" },
    @{ Cell = "E22"; Value = "  +8.45%  " },
    @{ Cell = "E23"; Value = "  +6.19%  " },
    @{ Cell = "D24"; Value = "11.26" },
    @{ Cell = "E24"; Value = "  +4.24%  " },
    @{ Cell = "D25"; Value = "2.284" },
    @{ Cell = "E25"; Value = "  +1.12%  " },
    @{ Cell = "D26"; Value = "2.169.97" },
    @{ Cell = "E26"; Value = "  +6.31%  " },
    @{ Cell = "D27"; Value = "160.67" },
    @{ Cell = "E27"; Value = "  +3.24%  " },
    @{ Cell = "D28"; Value = "20.21" },
    @{ Cell = "E28"; Value = "  +4.75%  " },
    @{ Cell = "D29"; Value = "2.200" },
    @{ Cell = "E29"; Value = "  +7.76%  " },
    @{ Cell = "D30"; Value = "5.650" },
    @{ Cell = "E30"; Value = "  +7.78%  " },
    @{ Cell = "D31"; Value = "122.33" },
    @{ Cell = "E31"; Value = "  +4.06%  " },
    @{ Cell = "D32"; Value = "1.028" },
    @{ Cell = "E32"; Value = "  +10.40%  " },
    @{ Cell = "D33"; Value = "0.09656" },
    @{ Cell = "E33"; Value = "  +3.96%  " },
    @{ Cell = "D34"; Value = "1.474" },
    @{ Cell = "E34"; Value = "  +12.32%  " },
    @{ Cell = "D35"; Value = "3.681" },
    @{ Cell = "E35"; Value = "  +3.21%  " },
    @{ Cell = "D36"; Value = "5.496" },
    @{ Cell = "E36"; Value = "  +5.47%  " },
    @{ Cell = "D37"; Value = "0.06300" },
    @{ Cell = "E37"; Value = "  +6.70%  " },
    @{ Cell = "E38"; Value = "  +6.82%  " },
    @{ Cell = "D39"; Value = "8.611" },
    @{ Cell = "E39"; Value = "  +6.39%  " },
    @{ Cell = "D40"; Value = "1.197" },
    @{ Cell = "E40"; Value = "  +5.12%  " },
    @{ Cell = "D41"; Value = "0.6110" },
    @{ Cell = "E41"; Value = "  +6.55%  " },
    @{ Cell = "E42"; Value = "  +8.01%  " },
    @{ Cell = "E43"; Value = "  +5.12%  " },
    @{ Cell = "D44"; Value = "1.001" },
    @{ Cell = "E44"; Value = "  -0.13%  " },
    @{ Cell = "D45"; Value = "1.272" },
    @{ Cell = "E45"; Value = "  +0.39%  " },
    @{ Cell = "D46"; Value = "12.62" },
    @{ Cell = "E46"; Value = "  +6.69%  " },
    @{ Cell = "B47"; Value = "Decentraland" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" },
    @{ Cell = "D47"; Value = "0.5726" },
    @{ Cell = "E47"; Value = "  +6.36%  " },
    @{ Cell = "B48"; Value = "RenderToken" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" },
    @{ Cell = "D48"; Value = "2.346" },
    @{ Cell = "E48"; Value = "  +31.39%  " },
    @{ Cell = "B49"; Value = "NEARProtocol" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" },
    @{ Cell = "D49"; Value = "2.001" },
    @{ Cell = "E49"; Value = "  +7.33%  " },
    @{ Cell = "B50"; Value = "Cronos" },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" },
    @{ Cell = "D50"; Value = "0.07421" },
    @{ Cell = "E50"; Value = "  +12.90%  " },
    @{ Cell = "D51"; Value = "114.42" },
    @{ Cell = "E51"; Value = "  +4.44%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $text = $u.Value

    # Several "Price" values (column D) are plain decimals such as
    # "341.74" or "0.4781". Assigning those strings straight to .Value
    # lets Excel's type-inference turn them into numeric cells (losing
    # trailing zeros / exact text and picking up float rounding noise),
    # whereas the source workbook stores every one of these as literal
    # text. Forcing a text number-format before the write keeps the
    # cell's type as text; clearing the format afterwards removes the
    # now-unneeded style override so the cell's style stays at its
    # original (default) index.
    if ($text -match '^-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}
